# Apply the changes described by the commit:
#  - rename Sheet1 -> "services", Sheet2 -> "other ports" (Sheet3 unchanged)
#  - populate the (until now empty) "other ports" sheet with a small table
#    that mirrors the "caGrid 1.0 Training" block already present on
#    "services" (A32:D34), re-using its formatting
#  - restore the various cursor/selection bookkeeping touched by the edit

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- rename the sheets -----------------------------------------------------
$ws1.Name = "services"
$ws2.Name = "other ports"

# --- new table on "other ports" --------------------------------------------
# Set the values first (and in this particular order) so that the shared
# string table is rebuilt with the same new-entry ordering as the original
# commit: usage.cagrid.org:55555, introduce, Protocol, UDP, Description,
# "introduce stats collecting".
$ws2.Range("A1").Value = "caGrid 1.0 Training"

$ws2.Range("A2").Value = "Host:Port"
$ws2.Range("B2").Value = "user"

$ws2.Range("A3").Value = "usage.cagrid.org:55555"
$ws2.Range("B3").Value = "introduce"

$ws2.Range("C2").Value = "Protocol"
$ws2.Range("C3").Value = "UDP"

$ws2.Range("D2").Value = "Description"
$ws2.Range("D3").Value = "introduce stats collecting"

# Re-use the formatting already used for the equivalent "caGrid 1.0 Training"
# block on the "services" sheet (title row A32, header+data rows A33:D34).
$ws1.Range("A32").Copy()
$ws2.Range("A1").PasteSpecial(-4122)   # xlPasteFormats

$ws1.Range("A33:D34").Copy()
$ws2.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

# Column widths for the new table.
$ws2.Columns.Item(1).ColumnWidth = 29.8
$ws2.Columns.Item(2).ColumnWidth = 22.65
$ws2.Columns.Item(3).ColumnWidth = 20.0
$ws2.Columns.Item(4).ColumnWidth = 53.5

# Row heights matching the title/header/data rows being mirrored.
$ws2.Rows.Item(1).RowHeight = 16.5
$ws2.Rows.Item(2).RowHeight = 17.25
$ws2.Rows.Item(3).RowHeight = 14.25

# --- selections --------------------------------------------------------
# "other ports" keeps D3 selected ...
$ws2.Range("D3").Select() | Out-Null

# ... while "services" stays the active tab, with its cursor moved to A46.
$ws1.Activate()
$ws1.Range("A46").Select() | Out-Null
